# Generate Report for Handoff
# The localization run finished and the handoff report needs to reflect the
# new status ("Ready for handoff") and the timestamps at which the handoff
# XLIFF files were (re)generated.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# Column E = zh-cn status, Column F = de-de status, Column G = latest
# handoff xliff-generation datetime for the (shared) de-de / overview entry.
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$ovw.Range("G2").Value = "2016-08-31 10:43:54"

# --- zh-cn sheet ------------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-08-31 10:43:50"

# --- de-de sheet ------------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$de.Range("C2").Value = "Ready for handoff"
$de.Range("H2").Value = "2016-08-31 10:43:54"

# --- Cosmetic: the Status/Latest-Handoff-Datetime columns grew a bit wider
# once the longer "Ready for handoff" text was in place (column autofit).
$ovw.Columns.Item(5).ColumnWidth = 16.3827
$ovw.Columns.Item(6).ColumnWidth = 16.3827
$zh.Columns.Item(3).ColumnWidth  = 16.3827
$de.Columns.Item(3).ColumnWidth  = 16.3827
